# Scheduled-runner update: refresh cached market-board figures
# (currentAveragePrice / NQ / HQ / LevePrice / LeveProfit columns)
# for a handful of leves across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 119086.84
$ws.Range("I40").Value = 752220
$ws.Range("J40").Value = 3971.7273
$ws.Range("K40").Value = 752220
$ws.Range("L40").Value = 3971.7273
$ws.Range("M40").Value = -752045
$ws.Range("N40").Value = -4321.7273

$ws.Range("H64").Value = 9500
$ws.Range("J64").Value = 9500
$ws.Range("L64").Value = 9500
$ws.Range("N64").Value = -9996

$ws.Range("H67").Value = 9500
$ws.Range("J67").Value = 9500
$ws.Range("L67").Value = 9500
$ws.Range("N67").Value = -11216

$ws.Range("H74").Value = 6199.7
$ws.Range("I74").Value = 5166.3335
$ws.Range("K74").Value = 5166.3335
$ws.Range("M74").Value = -4230.3335

$ws.Range("H77").Value = 6199.7
$ws.Range("I77").Value = 5166.3335
$ws.Range("K77").Value = 25831.6675
$ws.Range("M77").Value = -21151.6675

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws.Range("H138").Value = 2488.4
$ws.Range("I138").Value = 2129.1428
$ws.Range("K138").Value = 6387.428400000001
$ws.Range("M138").Value = -1247.428400000001

$ws.Range("H139").Value = 50000
$ws.Range("J139").Value = 50000
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280

$ws.Range("H141").Value = 4715.5557
$ws.Range("I141").Value = 4492.2856
$ws.Range("J141").Value = 5497
$ws.Range("K141").Value = 13476.8568
$ws.Range("L141").Value = 16491
$ws.Range("M141").Value = -8296.856800000001
$ws.Range("N141").Value = -26851

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1244.9
$ws.Range("I2").Value = 1200.4667
$ws.Range("K2").Value = 1200.4667
$ws.Range("M2").Value = -1087.4667

$ws.Range("H32").Value = 10686629
$ws.Range("I32").Value = 5377904
$ws.Range("K32").Value = 5377904
$ws.Range("M32").Value = -5377617

$ws.Range("H45").Value = 1920.1111
$ws.Range("I45").Value = 1035.25
$ws.Range("K45").Value = 1035.25
$ws.Range("M45").Value = -658.25

$ws.Range("H61").Value = 2358.4482
$ws.Range("I61").Value = 2081.4167
$ws.Range("K61").Value = 2081.4167
$ws.Range("M61").Value = -1869.4167

$ws.Range("H63").Value = 4100
$ws.Range("I63").Value = 2800
$ws.Range("J63").Value = 5400
$ws.Range("K63").Value = 2800
$ws.Range("L63").Value = 5400
$ws.Range("M63").Value = -2114
$ws.Range("N63").Value = -6772

$ws.Range("H66").Value = 4100
$ws.Range("I66").Value = 2800
$ws.Range("J66").Value = 5400
$ws.Range("K66").Value = 14000
$ws.Range("L66").Value = 27000
$ws.Range("M66").Value = -10568
$ws.Range("N66").Value = -33864

$ws.Range("H116").Value = 1244.9
$ws.Range("I116").Value = 1200.4667
$ws.Range("K116").Value = 1200.4667
$ws.Range("M116").Value = 1093.5333

$ws.Range("H132").Value = 2690.205
$ws.Range("I132").Value = 2108.3572
$ws.Range("K132").Value = 6325.071599999999
$ws.Range("M132").Value = -3795.071599999999

$ws.Range("H136").Value = 2358.4482
$ws.Range("I136").Value = 2081.4167
$ws.Range("K136").Value = 6244.250100000001
$ws.Range("M136").Value = -3694.250100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1244.9
$ws.Range("I3").Value = 1200.4667
$ws.Range("K3").Value = 1200.4667
$ws.Range("M3").Value = -1086.4667

$ws.Range("H20").Value = 19765.285
$ws.Range("I20").Value = 22526.334
$ws.Range("K20").Value = 22526.334
$ws.Range("M20").Value = -22279.334

$ws.Range("H40").Value = 37488.5
$ws.Range("J40").Value = 37488.5
$ws.Range("L40").Value = 37488.5
$ws.Range("N40").Value = -38018.5

$ws.Range("H96").Value = 74999
$ws.Range("I96").Value = 2000
$ws.Range("J96").Value = 111498.5
$ws.Range("K96").Value = 2000
$ws.Range("L96").Value = 111498.5
$ws.Range("M96").Value = 746
$ws.Range("N96").Value = -116990.5

$ws.Range("H134").Value = 11908008
$ws.Range("I134").Value = 2383128
$ws.Range("K134").Value = 7149384
$ws.Range("M134").Value = -7146849

$ws.Range("H138").Value = 100000
$ws.Range("I138").Value = 80000
$ws.Range("J138").Value = 120000
$ws.Range("K138").Value = 80000
$ws.Range("L138").Value = 120000
$ws.Range("M138").Value = -74860
$ws.Range("N138").Value = -130280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 151.6
$ws.Range("I7").Value = 122.045456
$ws.Range("J7").Value = 368.33334
$ws.Range("K7").Value = 122.045456
$ws.Range("L7").Value = 368.33334
$ws.Range("M7").Value = -9.045456000000001
$ws.Range("N7").Value = -594.33334

$ws.Range("H31").Value = 3191.6667
$ws.Range("I31").Value = 2128.25
$ws.Range("K31").Value = 2128.25
$ws.Range("M31").Value = -1833.25

$ws.Range("H34").Value = 3191.6667
$ws.Range("I34").Value = 2128.25
$ws.Range("K34").Value = 2128.25
$ws.Range("M34").Value = -1926.25

$ws.Range("H132").Value = 2100.9412
$ws.Range("I132").Value = 1937.9
$ws.Range("J132").Value = 2333.8572
$ws.Range("K132").Value = 5813.700000000001
$ws.Range("L132").Value = 7001.571599999999
$ws.Range("M132").Value = -3283.700000000001
$ws.Range("N132").Value = -12061.5716

$ws.Range("H137").Value = 120001
$ws.Range("J137").Value = 120001
$ws.Range("L137").Value = 120001
$ws.Range("N137").Value = -130201

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1017.0769
$ws.Range("I5").Value = 632.3333
$ws.Range("K5").Value = 1896.9999
$ws.Range("M5").Value = -1784.9999

$ws.Range("H7").Value = 210
$ws.Range("I7").Value = 150
$ws.Range("K7").Value = 450
$ws.Range("M7").Value = -338

$ws.Range("H50").Value = 3974.6667
$ws.Range("I50").Value = 4062
$ws.Range("J50").Value = 3800
$ws.Range("K50").Value = 12186
$ws.Range("L50").Value = 11400
$ws.Range("M50").Value = -11705
$ws.Range("N50").Value = -12362

$ws.Range("H53").Value = 3974.6667
$ws.Range("I53").Value = 4062
$ws.Range("J53").Value = 3800
$ws.Range("K53").Value = 12186
$ws.Range("L53").Value = 11400
$ws.Range("M53").Value = -11705
$ws.Range("N53").Value = -12362

$ws.Range("H104").Value = 8343
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 8343
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 25029
$ws.Range("M104").ClearContents()
$ws.Range("N104").Value = -30271

$ws.Range("H118").Value = 1929
$ws.Range("I118").Value = 1929
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 5787
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -4544
$ws.Range("N118").ClearContents()

$ws.Range("H132").Value = 1843.28
$ws.Range("I132").Value = 1321.25
$ws.Range("J132").Value = 2088.9412
$ws.Range("K132").Value = 11891.25
$ws.Range("L132").Value = 18800.4708
$ws.Range("M132").Value = -9361.25
$ws.Range("N132").Value = -23860.4708

$ws.Range("H135").Value = 1017.0769
$ws.Range("I135").Value = 632.3333
$ws.Range("K135").Value = 5690.9997
$ws.Range("M135").Value = -3155.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5289.0386
$ws.Range("I132").Value = 4439.1113
$ws.Range("J132").Value = 7201.375
$ws.Range("K132").Value = 13317.3339
$ws.Range("L132").Value = 21604.125
$ws.Range("M132").Value = -10787.3339
$ws.Range("N132").Value = -26664.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6247.4062
$ws.Range("I40").Value = 5969.6553
$ws.Range("K40").Value = 5969.6553
$ws.Range("M40").Value = -5833.6553

$ws.Range("H46").Value = 3116.923
$ws.Range("I46").Value = 598
$ws.Range("J46").Value = 3872.6
$ws.Range("K46").Value = 598
$ws.Range("L46").Value = 3872.6
$ws.Range("M46").Value = -410
$ws.Range("N46").Value = -4248.6

$ws.Range("H100").Value = 3075.5
$ws.Range("I100").Value = 2860.8
$ws.Range("J100").Value = 3433.3333
$ws.Range("K100").Value = 2860.8
$ws.Range("L100").Value = 3433.3333
$ws.Range("M100").Value = -2319.8
$ws.Range("N100").Value = -4515.3333

$ws.Range("H132").Value = 4459.4
$ws.Range("I132").Value = 2919.6
$ws.Range("K132").Value = 8758.799999999999
$ws.Range("M132").Value = -6228.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1923.8422
$ws.Range("I100").Value = 1915.8125
$ws.Range("K100").Value = 3831.625
$ws.Range("M100").Value = -3290.625

$ws.Range("H132").Value = 3907.5
$ws.Range("I132").Value = 3409.5264
$ws.Range("K132").Value = 10228.5792
$ws.Range("M132").Value = -7698.5792
